$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.203.91"
$ws.Range("E2").Value = "  -2.79%  "
$ws.Range("D3").Value = "2.192.13"
$ws.Range("E3").Value = "  -7.44%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'296.66"
$ws.Range("E5").Value = "  -4.44%  "
$ws.Range("D6").Value = "'81.46"
$ws.Range("E6").Value = "  -5.34%  "
$ws.Range("D7").Value = "'0.508"
$ws.Range("E7").Value = "  -4.80%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.464"
$ws.Range("E9").Value = "  -5.28%  "
$ws.Range("D10").Value = "'0.0767"
$ws.Range("E10").Value = "  -7.37%  "
$ws.Range("D11").Value = "'28.88"
$ws.Range("E11").Value = "  -4.83%  "
$ws.Range("D12").Value = "'47.31"
$ws.Range("E13").Value = "  -2.79%  "
$ws.Range("D14").Value = "2.540.18"
$ws.Range("E14").Value = "  -7.22%  "
$ws.Range("D15").Value = "'6.23"
$ws.Range("E15").Value = "  -3.80%  "
$ws.Range("D16").Value = "'13.90"
$ws.Range("E16").Value = "  -7.17%  "
$ws.Range("D17").Value = "2.198.67"
$ws.Range("E17").Value = "  -7.30%  "
$ws.Range("D18").Value = "'0.707"
$ws.Range("E18").Value = "  -6.64%  "
$ws.Range("D19").Value = "39.117.66"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("D20").Value = "0.0₃0867"
$ws.Range("E20").Value = "  -4.63%  "
$ws.Range("D21").Value = "'5.69"
$ws.Range("E21").Value = "  -7.13%  "
$ws.Range("D22").Value = "'64.82"
$ws.Range("E22").Value = "  -5.03%  "
$ws.Range("D23").Value = "'10.21"
$ws.Range("E23").Value = "  -5.45%  "
$ws.Range("D24").Value = "'224.93"
$ws.Range("E24").Value = "  -4.27%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'2.39"
$ws.Range("E26").Value = "  -7.21%  "
$ws.Range("D27").Value = "'1.78"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").Value = "'22.42"
$ws.Range("E28").Value = "  -5.39%  "
$ws.Range("D29").Value = "'2.16"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").Value = "'9.02"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("D31").Value = "'148.59"
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("D32").Value = "'31.45"
$ws.Range("E32").Value = "  -8.16%  "
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").Value = "'4.80"
$ws.Range("E34").Value = "  -7.79%  "
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("D36").Value = "'0.0688"
$ws.Range("E36").Value = "  -5.89%  "
$ws.Range("E37").Value = "  -4.52%  "
$ws.Range("D38").Value = "'0.0963"
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("D39").Value = "'15.10"
$ws.Range("E39").Value = "  -5.41%  "
$ws.Range("E40").Value = "  -6.63%  "
$ws.Range("D41").Value = "'1.62"
$ws.Range("E41").Value = "  -5.14%  "
$ws.Range("D42").Value = "'3.60"
$ws.Range("E42").Value = "  -6.00%  "
$ws.Range("D43").Value = "1.889.89"
$ws.Range("E43").Value = "  -3.96%  "
$ws.Range("D44").Value = "'2.09"
$ws.Range("E44").Value = "  -12.19%  "
$ws.Range("E45").Value = "  -4.11%  "
$ws.Range("D46").Value = "'8.92"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D47").Value = "'15.92"
$ws.Range("E47").Value = "  -10.62%  "
$ws.Range("D48").Value = "'2.59"
$ws.Range("E48").Value = "  -3.56%  "
$ws.Range("D49").Value = "2.412.84"
$ws.Range("E49").Value = "  -7.02%  "
$ws.Range("D50").Value = "'71.19"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").Value = "'86.48"
$ws.Range("E51").Value = "  -7.10%  "
